$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at 457, shifting existing rows 457:505 down to 458:506.
$ws.Rows("457:457").Insert()

# Populate the newly inserted row 457 with the new weekly data point.
$ws.Cells.Item(457, 1).Value = 4
$ws.Cells.Item(457, 2).Value = "Feria Lagunitas de Puerto Montt"
$ws.Cells.Item(457, 3).Value = "Los Lagos"
$ws.Cells.Item(457, 4).Value = 45194
$ws.Cells.Item(457, 5).Value = 10
$ws.Cells.Item(457, 6).Value = 100112043
$ws.Cells.Item(457, 7).Value = "Pepino ensalada"
$ws.Cells.Item(457, 8).Value = "Sin especificar"
$ws.Cells.Item(457, 9).Value = "Primera"
$ws.Cells.Item(457, 10).Value = 80
$ws.Cells.Item(457, 11).Value = 17000
$ws.Cells.Item(457, 12).Value = 17000
$ws.Cells.Item(457, 13).Value = 17000
$ws.Cells.Item(457, 14).Value = "$/caja 60 unidades"
$ws.Cells.Item(457, 15).Value = "Región de Arica y Parinacota"
$ws.Cells.Item(457, 16).Value = 283
$ws.Cells.Item(457, 17).Value = 60
$ws.Cells.Item(457, 18).Value = "Hortaliza"
